# controle_processos.xlsx - "alteração do módulo de controle de contratos"
#
# 1) Remove workbook protection marker
# 2) Update a handful of "etapa" (K column) values
# 3) Clear a bunch of stray empty cells (trailing blank columns from the
#    original CSV/export) that should not be present any more
# 4) Apply explicit column widths to the main data sheet
# 5) Normalize page margins to Excel's defaults

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Drop workbook protection (was an empty/no-op protection marker) ---
$wb.Unprotect()

# --- 2) Data updates on the "etapa" column ---
$ws.Range("K2").Value = "Em recurso"
$ws.Range("K3").Value = "Em recurso"
$ws.Range("K4").Value = "Em recurso"
$ws.Range("K7").Value = "CJACM"

# --- 3) Clear stray empty cells that trailed several rows ---
$emptyCells = @(
    "L3", "L4", "L5",
    "L12", "L13", "L14", "L15", "L16", "L17", "L18", "L19",
    "L20", "L21", "L22", "L23", "L24", "L25", "L26", "L27", "L28", "L29",
    "L30", "L31", "L32", "L33", "L34", "L35", "L36", "L37", "L38", "L39",
    "L40", "L41", "L42", "L43", "L44", "L45", "L46", "L47", "L48", "L49",
    "L50", "L51", "L52", "L53", "L54", "L55", "L56", "L57", "L58", "L59",
    "L60", "L61", "L62", "L63",
    "L66",
    "L74", "L75", "L76",
    "E76", "I76", "J76"
)
foreach ($addr in $emptyCells) {
    $ws.Range($addr).ClearContents()
}

# --- 4) Explicit column widths (values compensate for Excel's standard
#         +0.83 "internal padding" so the saved width matches exactly) ---
$ws.Columns.Item(1).ColumnWidth = 9.1666666   # A
$ws.Columns.Item(2).ColumnWidth = 9.1666666   # B
$ws.Columns.Item(3).ColumnWidth = 24.1666666  # C
$ws.Columns.Item(4).ColumnWidth = 34.1666666  # D
$ws.Columns.Item(6).ColumnWidth = 39.1666666  # F
$ws.Columns.Item(7).ColumnWidth = 9.1666666   # G
$ws.Columns.Item(8).ColumnWidth = 19.1666666  # H
$ws.Columns.Item(9).ColumnWidth = 9.1666666   # I
$ws.Columns.Item(10).ColumnWidth = 19.1666666 # J
$ws.Columns.Item(11).ColumnWidth = 19.1666666 # K

# --- 5) Page margins back to Excel's standard defaults ---
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72

# Make sure the sheet is the selected/active tab on reopen.
$ws.Select()
$ws.Range("A1").Select()
